# Aggiornamento fino al 21 aprile: aggiunta delle righe 230-233
# (dati per le date 2021-04-18 .. 2021-04-21) in coda al foglio "Sheet1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# La colonna A usa lo stile "data" applicato alla riga precedente (229):
# copiamo la formattazione dell'ultima riga esistente (A229:D229) sulle
# nuove righe, cosi' la colonna A eredita lo stesso formato data (s="2")
# mentre B:D restano senza formattazione specifica, come nel resto del foglio.
$ws.Range("A229:D229").Copy()
$ws.Range("A230:D233").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @(230, 44304, 1, 6, 68.99724011039559),
    @(231, 44305, 0, 4, 45.99816007359706),
    @(232, 44306, 1, 5, 57.49770009199631),
    @(233, 44307, 0, 5, 57.49770009199631)
)

foreach ($entry in $newRows) {
    $r        = $entry[0]
    $dataSer  = $entry[1]
    $nuoviPos = $entry[2]
    $somma7gg = $entry[3]
    $somma100k = $entry[4]

    $ws.Cells.Item($r, 1).Value = $dataSer
    $ws.Cells.Item($r, 2).Value = $nuoviPos
    $ws.Cells.Item($r, 3).Value = $somma7gg
    $ws.Cells.Item($r, 4).Value = $somma100k
}
